$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new row of data (row 3)
$ws.Range("A3").Value = "ايمن محمد بيومي "
$ws.Range("B3").Value = "aymanmohamed@gmai.com"
$ws.Range("C3").Value = "https://github.com/aymanmohamed78/Security-Task.git"

# Match the existing look of row 2 (border/fill style) for the new data cells
$ws.Range("B2:C2").Copy()
$ws.Range("B3:C3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Wire up the hyperlinks for the new email and repo link cells
$ws.Hyperlinks.Add($ws.Range("B3"), "mailto:aymanmohamed@gmai.com")
$ws.Hyperlinks.Add($ws.Range("C3"), "https://github.com/aymanmohamed78/Security-Task.git")

# Re-apply the shared style after adding hyperlinks (Hyperlinks.Add can restyle the cell)
$ws.Range("B2:C2").Copy()
$ws.Range("B3:C3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Move selection to C3 to mirror the saved workbook state
[void]$ws.Range("C3").Select()
